$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the eval-data JSON stored in column D (shared string) ---
# The original text included a "completionTokens" field that has been
# dropped from the eval payload; every row in D2:D51 shares this text.
$newText = '{"fcCount":1,"fcInfo":{"apiair-conditionerupdate_POST":16}}'
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 4).Value = $newText
}

# --- Update the active sheet's view / selection ---
$ws.Activate() | Out-Null
$ws.Range("E2:N68").Select() | Out-Null
